# Applies the "Finished EEMs prep for Alan" edit to the TMP_FEOM_CO_1 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Mark a batch of EEM rows as "Omit" in the Action column (C).
$omitRows = @(25,26,27,28,29,30,31,35,36,37,38,39,40,41,42,43,47,48)
foreach ($r in $omitRows) {
    $ws.Range("C$r").Value = "Omit"
}

# 2) Fix a mislabeled sample id: row 65 was a duplicate "AW2.01.C" label;
#    it should read "AW3.01.C".
$ws.Range("A65").Value = "AW3.01.C"

# 3) Remove the stray duplicate "AW2.ASW" row (old row 92) - everything
#    below it shifts up by one row.
$ws.Rows(92).Delete()

# Deleting that row breaks the formula chain for the row that slides into
# position 92 (it referenced the now-deleted row), so replace it with the
# literal, already-computed value - matching what the chain would have
# produced.
$ws.Range("B92").Value = 8

# 4) Restore the view: scroll so row 82 is at the top and select C94.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
[void]$ws.Range("C94").Select()
